$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.519.63"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.178.17"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.88"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.79"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.175.24"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("E11").Value = "  +5.71%  "
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.19"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.707.09"
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.247.12"
$ws.Range("E17").Value = "  +4.49%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.38"
$ws.Range("E18").Value = "  +3.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.325.31"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.93"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.752"
$ws.Range("E22").Value = "  +3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.79"
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.49"
$ws.Range("E24").Value = "  +11.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.70"
$ws.Range("E25").Value = "  +5.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.28"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").Value = "  +7.89%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.75"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.34"
$ws.Range("E32").Value = "  +4.60%  "
$ws.Range("E33").Value = "  +11.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.52"
$ws.Range("E34").Value = "  +6.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("E37").Value = "  +4.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "472.10"
$ws.Range("E40").Value = "  +6.85%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.90"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.43"
$ws.Range("E42").Value = "  +8.41%  "
$ws.Range("E43").Value = "  +8.76%  "
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.945.15"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.35"
$ws.Range("E47").Value = "  +6.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.55"
$ws.Range("E48").Value = "  +5.20%  "
$ws.Range("E49").Value = "  +6.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.14"
$ws.Range("E51").Value = "  +3.66%  "
